$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.590.20'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '2.645.16'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.10%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.40%  '

$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("E11").Value = '  +4.68%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("D14").Value = '3.118.38'
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").Value = '63.385.57'
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("E16").Value = '  +2.03%  '

$ws.Range("D17").Value = '2.644.04'
$ws.Range("E17").Value = '  -0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.22%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.77%  '

$ws.Range("E25").Value = '  +2.04%  '

$ws.Range("B26").Value = 'Bittensor'
$ws.Range("C26").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '582.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.50%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.48%  '

$ws.Range("E29").Value = '  -1.41%  '

$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.73%  '

$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("E32").Value = '  +4.12%  '

$ws.Range("E33").Value = '  -2.38%  '

$ws.Range("D34").Value = '0.0₃0829'
$ws.Range("E34").Value = '  +3.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '166.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.407'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.47%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '168.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0571'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.631'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0247'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0963'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.178'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.73%  '
